$wb = $excel.ActiveWorkbook

# --- Rename sheets: unify DataNode / DataTable naming convention ---
$wb.Worksheets("Property1").Name = "DataNode_1"
$wb.Worksheets("Property2").Name = "DataNode_2"
$wb.Worksheets("Record_Hero").Name = "DataTable_Hero"
$wb.Worksheets("Record_Bag").Name = "DataTable_Bag"
$wb.Worksheets("Record_CommPropertyValue").Name = "DataTable_CommPropertyValue"
$wb.Worksheets("Record_Task").Name = "DataTable_Task"
# "Component" sheet name stays the same

# --- Update the stored selection on DataTable_Task before it loses focus ---
$wb.Worksheets("DataTable_Task").Range("I43").Select() | Out-Null

# --- Update the stored selection on Component before it loses focus ---
$wb.Worksheets("Component").Range("K40").Select() | Out-Null

# --- Remove the obsolete Record_Building sheet entirely ---
$excel.DisplayAlerts = $false
$wb.Worksheets("Record_Building").Delete() | Out-Null
$excel.DisplayAlerts = $true

# --- Make DataTable_Hero the active tab (was Property1/DataNode_1 before) ---
$wb.Worksheets("DataTable_Hero").Activate()
